$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 802
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 802
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 802
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3834.423
$ws.Range("I116").Value = 3022.9167
$ws.Range("J116").Value = 4530
$ws.Range("K116").Value = 3022.9167
$ws.Range("L116").Value = 4530
$ws.Range("M116").Value = 419.0832999999998
$ws.Range("N116").Value = -11414

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1768.0834
$ws.Range("I110").Value = 613.2308
$ws.Range("J110").Value = 3132.9092
$ws.Range("K110").Value = 613.2308
$ws.Range("L110").Value = 3132.9092
$ws.Range("M110").Value = 1431.7692
$ws.Range("N110").Value = -7222.9092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2051.5217
$ws.Range("I20").Value = 1683.5385
$ws.Range("J20").Value = 2529.9
$ws.Range("K20").Value = 1683.5385
$ws.Range("L20").Value = 2529.9
$ws.Range("M20").Value = -1436.5385
$ws.Range("N20").Value = -3023.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1045957.5
$ws.Range("I86").Value = 1233519.9
$ws.Range("J86").Value = 170666.67
$ws.Range("K86").Value = 1233519.9
$ws.Range("L86").Value = 170666.67
$ws.Range("M86").Value = -1232396.9
$ws.Range("N86").Value = -172912.67

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1045957.5
$ws.Range("I89").Value = 1233519.9
$ws.Range("J89").Value = 170666.67
$ws.Range("K89").Value = 6167599.5
$ws.Range("L89").Value = 853333.3500000001
$ws.Range("M89").Value = -6161983.5
$ws.Range("N89").Value = -864565.3500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1397.2222
$ws.Range("I105").Value = 1302.9412
$ws.Range("K105").Value = 1302.9412
$ws.Range("M105").Value = 444.0588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2109
$ws.Range("I107").Value = 1085.1666
$ws.Range("K107").Value = 1085.1666
$ws.Range("M107").Value = 834.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1002526
$ws.Range("I31").Value = 1924768.2
$ws.Range("J31").Value = 3430.1667
$ws.Range("K31").Value = 1924768.2
$ws.Range("L31").Value = 3430.1667
$ws.Range("M31").Value = -1924473.2
$ws.Range("N31").Value = -4020.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1002526
$ws.Range("I34").Value = 1924768.2
$ws.Range("J34").Value = 3430.1667
$ws.Range("K34").Value = 1924768.2
$ws.Range("L34").Value = 3430.1667
$ws.Range("M34").Value = -1924566.2
$ws.Range("N34").Value = -3834.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9275.333000000001
$ws.Range("J51").Value = 10409.75
$ws.Range("L51").Value = 10409.75
$ws.Range("N51").Value = -11881.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11366916
$ws.Range("I58").Value = 2144.9565
$ws.Range("J58").Value = 23814046
$ws.Range("K58").Value = 2144.9565
$ws.Range("L58").Value = 23814046
$ws.Range("M58").Value = -1941.9565
$ws.Range("N58").Value = -23814452

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 9275.333000000001
$ws.Range("J61").Value = 10409.75
$ws.Range("L61").Value = 10409.75
$ws.Range("N61").Value = -11105.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2593.6667
$ws.Range("I94").Value = 2431.7144
$ws.Range("J94").Value = 2660.353
$ws.Range("K94").Value = 2431.7144
$ws.Range("L94").Value = 2660.353
$ws.Range("M94").Value = -1980.7144
$ws.Range("N94").Value = -3562.353

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2070.125
$ws.Range("I122").Value = 1641.6
$ws.Range("J122").Value = 3600.5715
$ws.Range("K122").Value = 4924.799999999999
$ws.Range("L122").Value = 10801.7145
$ws.Range("M122").Value = -2474.799999999999
$ws.Range("N122").Value = -15701.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2432.6667
$ws.Range("I132").Value = 2163
$ws.Range("J132").Value = 2972
$ws.Range("K132").Value = 6489
$ws.Range("L132").Value = 8916
$ws.Range("M132").Value = -3959
$ws.Range("N132").Value = -13976

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 11366916
$ws.Range("I136").Value = 2144.9565
$ws.Range("J136").Value = 23814046
$ws.Range("K136").Value = 6434.869499999999
$ws.Range("L136").Value = 71442138
$ws.Range("M136").Value = -3884.869499999999
$ws.Range("N136").Value = -71447238

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 438.5625
$ws.Range("I12").Value = 20.454546
$ws.Range("J12").Value = 657.5714
$ws.Range("K12").Value = 61.363638
$ws.Range("L12").Value = 1972.7142
$ws.Range("M12").Value = 111.636362
$ws.Range("N12").Value = -2318.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1164.3131
$ws.Range("I131").Value = 3072.5
$ws.Range("J131").Value = 996.5604
$ws.Range("K131").Value = 9217.5
$ws.Range("L131").Value = 2989.6812
$ws.Range("M131").Value = -4177.5
$ws.Range("N131").Value = -13069.6812

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3520.9092
$ws.Range("I80").Value = 2956
$ws.Range("J80").Value = 3991.6667
$ws.Range("K80").Value = 2956
$ws.Range("L80").Value = 3991.6667
$ws.Range("M80").Value = -1958
$ws.Range("N80").Value = -5987.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3520.9092
$ws.Range("I83").Value = 2956
$ws.Range("J83").Value = 3991.6667
$ws.Range("K83").Value = 14780
$ws.Range("L83").Value = 19958.3335
$ws.Range("M83").Value = -9788
$ws.Range("N83").Value = -29942.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 51351.332
$ws.Range("I102").Value = 2832.6
$ws.Range("J102").Value = 172648.17
$ws.Range("K102").Value = 2832.6
$ws.Range("L102").Value = 172648.17
$ws.Range("M102").Value = -1210.6
$ws.Range("N102").Value = -175892.17

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3539.6545
$ws.Range("I132").Value = 3491.4324
$ws.Range("K132").Value = 10474.2972
$ws.Range("M132").Value = -7944.297200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 52633068
$ws.Range("I16").Value = 66667724
$ws.Range("J16").Value = 3095.5
$ws.Range("K16").Value = 66667724
$ws.Range("L16").Value = 3095.5
$ws.Range("M16").Value = -66667554
$ws.Range("N16").Value = -3435.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1514.2222
$ws.Range("I22").Value = 427.2
$ws.Range("J22").Value = 2873
$ws.Range("K22").Value = 427.2
$ws.Range("L22").Value = 2873
$ws.Range("M22").Value = -132.2
$ws.Range("N22").Value = -3463

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1514.2222
$ws.Range("I27").Value = 427.2
$ws.Range("J27").Value = 2873
$ws.Range("K27").Value = 427.2
$ws.Range("L27").Value = 2873
$ws.Range("M27").Value = -320.2
$ws.Range("N27").Value = -3087

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2888.8635
$ws.Range("I122").Value = 2453.8
$ws.Range("J122").Value = 3821.1428
$ws.Range("K122").Value = 7361.400000000001
$ws.Range("L122").Value = 11463.4284
$ws.Range("M122").Value = -4911.400000000001
$ws.Range("N122").Value = -16363.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2657.4688
$ws.Range("I136").Value = 3172.8462
$ws.Range("J136").Value = 2304.842
$ws.Range("K136").Value = 9518.5386
$ws.Range("L136").Value = 6914.526
$ws.Range("M136").Value = -6968.5386
$ws.Range("N136").Value = -12014.526
